$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0864294738795979
$ws.Range("C2").Value = 0.1717255535034242
$ws.Range("D2").Value = 0.05132172044861762
$ws.Range("E2").Value = 0.2265429770454551
$ws.Range("F2").Value = 0.2173128427139068

$ws.Range("B3").Value = 0.1015794586745889
$ws.Range("C3").Value = 0.1777261608473147
$ws.Range("D3").Value = 0.06528206036901001
$ws.Range("E3").Value = 0.255503542771935
$ws.Range("F3").Value = 0.244016352129155

$ws.Range("B4").Value = 0.09564297872977319
$ws.Range("C4").Value = 0.182775340538003
$ws.Range("D4").Value = 0.07239544141620145
$ws.Range("E4").Value = 0.2690640098864979
$ws.Range("F4").Value = 0.2626740711899914

$ws.Range("B5").Value = 0.0515887520981297
$ws.Range("C5").Value = 0.1373511933772956
$ws.Range("D5").Value = 0.03807726663693461
$ws.Range("E5").Value = 0.1951339709966838
$ws.Range("F5").Value = 0.1973764272229122

$ws.Range("B6").Value = 0.03412360743610758
$ws.Range("C6").Value = 0.1315810182071815
$ws.Range("D6").Value = 0.02988307038476774
$ws.Range("E6").Value = 0.1728672044801088
$ws.Range("F6").Value = 0.1786326142932414

$ws.Range("B7").Value = 0.03856004825347423
$ws.Range("C7").Value = 0.1426784820072856
$ws.Range("D7").Value = 0.04033731687300466
$ws.Range("E7").Value = 0.2008415217852241
$ws.Range("F7").Value = 0.2090615806303401

$ws.Range("B8").Value = 0.06589614311488577
$ws.Range("C8").Value = 0.08267801794971415
$ws.Range("D8").Value = 0.01154809474451304
$ws.Range("E8").Value = 0.107462061884709
$ws.Range("F8").Value = 0.09298898687755791

$ws.Range("B9").Value = -0.03841838264255042
$ws.Range("C9").Value = 0.03841838264255042
$ws.Range("D9").Value = 0.001596828758587614
$ws.Range("E9").Value = 0.03996033982072242
$ws.Range("F9").Value = 0.01346420998712108

$ws.Range("B10").Value = -0.02233679963354765
$ws.Range("C10").Value = 0.02233679963354765
$ws.Range("D10").Value = 0.0004989326178692543
$ws.Range("E10").Value = 0.02233679963354765
